$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3   = 5.545
    14  = 5.755999999999999
    16  = 5.389
    21  = 8.996
    23  = 7.398999999999999
    25  = 6.396
    26  = 6.547
    29  = 5.689
    40  = 9.137
    53  = 5.748
    57  = 5.090999999999999
    59  = 4.435
    65  = 5.746
    69  = 5.095
    79  = 5.623
    83  = 5.523999999999999
    91  = 5.468999999999999
    93  = 5.635000000000001
    100 = 5.931
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 2).Value = $updates[$row]
}
